$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Log an additional 8 hours worked on the date in row 19 (D19)
$ws.Range("D19").Value = 8

# Update the active selection to H11, matching the author's last cursor position
$ws.Range("H11").Select()
